$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3052.7778
$ws.Range("I62").Value = 3244.5
$ws.Range("J62").Value = 2899.4
$ws.Range("K62").Value = 3244.5
$ws.Range("L62").Value = 2899.4
$ws.Range("M62").Value = -2620.5
$ws.Range("N62").Value = -4147.4

$ws.Range("H65").Value = 3052.7778
$ws.Range("I65").Value = 3244.5
$ws.Range("J65").Value = 2899.4
$ws.Range("K65").Value = 16222.5
$ws.Range("L65").Value = 14497
$ws.Range("M65").Value = -13102.5
$ws.Range("N65").Value = -20737

$ws.Range("H86").Value = 3110
$ws.Range("I86").Value = 2900
$ws.Range("J86").Value = 3460
$ws.Range("K86").Value = 2900
$ws.Range("L86").Value = 3460
$ws.Range("M86").Value = -1777
$ws.Range("N86").Value = -5706

$ws.Range("H89").Value = 3110
$ws.Range("I89").Value = 2900
$ws.Range("J89").Value = 3460
$ws.Range("K89").Value = 14500
$ws.Range("L89").Value = 17300
$ws.Range("M89").Value = -8884
$ws.Range("N89").Value = -28532

$ws.Range("H121").Value = 614.91174
$ws.Range("I121").Value = 1000
$ws.Range("J121").Value = 603.24243
$ws.Range("K121").Value = 3000
$ws.Range("L121").Value = 1809.72729
$ws.Range("M121").Value = -1253
$ws.Range("N121").Value = -5303.72729

$ws.Range("H137").Value = 1603.1464
$ws.Range("I137").Value = 1331.7632
$ws.Range("J137").Value = 5040.6665
$ws.Range("K137").Value = 3995.2896
$ws.Range("L137").Value = 15121.9995
$ws.Range("M137").Value = -1445.2896
$ws.Range("N137").Value = -20221.9995

$ws.Range("H138").Value = 3354.8384
$ws.Range("I138").Value = 2253.9343
$ws.Range("J138").Value = 5122.079
$ws.Range("K138").Value = 6761.8029
$ws.Range("L138").Value = 15366.237
$ws.Range("M138").Value = -1621.8029
$ws.Range("N138").Value = -25646.237

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 26500
$ws.Range("I3").Value = 26500
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 26500
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -26385
$ws.Range("N3").ClearContents()

$ws.Range("H39").Value = 10949.7
$ws.Range("I39").Value = 10166.333
$ws.Range("J39").Value = 18000
$ws.Range("K39").Value = 10166.333
$ws.Range("L39").Value = 18000
$ws.Range("M39").Value = -9646.333000000001
$ws.Range("N39").Value = -19040

$ws.Range("H61").Value = 2225.2068
$ws.Range("I61").Value = 2004.4783
$ws.Range("J61").Value = 3071.3333
$ws.Range("K61").Value = 2004.4783
$ws.Range("L61").Value = 3071.3333
$ws.Range("M61").Value = -1792.4783
$ws.Range("N61").Value = -3495.3333

$ws.Range("H74").Value = 1393.4193
$ws.Range("I74").Value = 1370.2593
$ws.Range("J74").Value = 1549.75
$ws.Range("K74").Value = 1370.2593
$ws.Range("L74").Value = 1549.75
$ws.Range("M74").Value = -496.2592999999999
$ws.Range("N74").Value = -3297.75

$ws.Range("H77").Value = 1393.4193
$ws.Range("I77").Value = 1370.2593
$ws.Range("J77").Value = 1549.75
$ws.Range("K77").Value = 6851.296499999999
$ws.Range("L77").Value = 7748.75
$ws.Range("M77").Value = -2483.296499999999
$ws.Range("N77").Value = -16484.75

$ws.Range("H110").Value = 922.9167
$ws.Range("I110").Value = 926.0909
$ws.Range("J110").Value = 888
$ws.Range("K110").Value = 926.0909
$ws.Range("L110").Value = 888
$ws.Range("M110").Value = 1118.9091
$ws.Range("N110").Value = -4978

$ws.Range("H132").Value = 1798.409
$ws.Range("I132").Value = 1487.6842
$ws.Range("J132").Value = 3766.3333
$ws.Range("K132").Value = 4463.0526
$ws.Range("L132").Value = 11298.9999
$ws.Range("M132").Value = -1933.0526
$ws.Range("N132").Value = -16358.9999

$ws.Range("H136").Value = 2225.2068
$ws.Range("I136").Value = 2004.4783
$ws.Range("J136").Value = 3071.3333
$ws.Range("K136").Value = 6013.4349
$ws.Range("L136").Value = 9213.999899999999
$ws.Range("M136").Value = -3463.4349
$ws.Range("N136").Value = -14313.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 45047.668
$ws.Range("I20").Value = 61026.65
$ws.Range("J20").Value = 6241.5713
$ws.Range("K20").Value = 61026.65
$ws.Range("L20").Value = 6241.5713
$ws.Range("M20").Value = -60779.65
$ws.Range("N20").Value = -6735.5713

$ws.Range("H134").Value = 3830.1072
$ws.Range("I134").Value = 2735.7
$ws.Range("J134").Value = 6566.125
$ws.Range("K134").Value = 8207.099999999999
$ws.Range("L134").Value = 19698.375
$ws.Range("M134").Value = -5672.099999999999
$ws.Range("N134").Value = -24768.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H132").Value = 1562.5264
$ws.Range("I132").Value = 1334.6177
$ws.Range("K132").Value = 4003.8531
$ws.Range("M132").Value = -1473.8531

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1579.5714
$ws.Range("I5").Value = 1527.2222
$ws.Range("J5").Value = 1673.8
$ws.Range("K5").Value = 4581.6666
$ws.Range("L5").Value = 5021.4
$ws.Range("M5").Value = -4469.6666
$ws.Range("N5").Value = -5245.4

$ws.Range("H131").Value = 876.23
$ws.Range("J131").Value = 888.13684
$ws.Range("L131").Value = 2664.41052
$ws.Range("N131").Value = -12744.41052

$ws.Range("H132").Value = 1891.8
$ws.Range("I132").Value = 1500.7273
$ws.Range("J132").Value = 2118.2104
$ws.Range("K132").Value = 13506.5457
$ws.Range("L132").Value = 19063.8936
$ws.Range("M132").Value = -10976.5457
$ws.Range("N132").Value = -24123.8936

$ws.Range("H133").Value = 3393.2222
$ws.Range("I133").Value = 1399.5
$ws.Range("K133").Value = 4198.5
$ws.Range("M133").Value = 861.5

$ws.Range("H134").Value = 4741.115
$ws.Range("I134").Value = 2630.6428
$ws.Range("J134").Value = 7203.3335
$ws.Range("K134").Value = 7891.928400000001
$ws.Range("L134").Value = 21610.0005
$ws.Range("M134").Value = -2821.928400000001
$ws.Range("N134").Value = -31750.0005

$ws.Range("H135").Value = 1579.5714
$ws.Range("I135").Value = 1527.2222
$ws.Range("J135").Value = 1673.8
$ws.Range("K135").Value = 13744.9998
$ws.Range("L135").Value = 15064.2
$ws.Range("M135").Value = -11209.9998
$ws.Range("N135").Value = -20134.2

$ws.Range("H136").Value = 4074.125
$ws.Range("J136").Value = 4991.75
$ws.Range("L136").Value = 14975.25
$ws.Range("N136").Value = -25175.25

$ws.Range("H137").Value = 47622760
$ws.Range("I137").Value = 3515
$ws.Range("J137").Value = 66670456
$ws.Range("K137").Value = 10545
$ws.Range("L137").Value = 200011368
$ws.Range("M137").Value = -5445
$ws.Range("N137").Value = -200021568

$ws.Range("H139").Value = 2219.7778
$ws.Range("I139").Value = 1824.1666
$ws.Range("J139").Value = 3011
$ws.Range("K139").Value = 5472.4998
$ws.Range("L139").Value = 9033
$ws.Range("M139").Value = -332.4997999999996
$ws.Range("N139").Value = -19313

$ws.Range("H140").Value = 2286.111
$ws.Range("I140").Value = 1567.875
$ws.Range("J140").Value = 8032
$ws.Range("K140").Value = 4703.625
$ws.Range("L140").Value = 24096
$ws.Range("M140").Value = 476.375
$ws.Range("N140").Value = -34456

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1688.4445
$ws.Range("I113").Value = 856
$ws.Range("J113").Value = 2008.6154
$ws.Range("K113").Value = 856
$ws.Range("L113").Value = 2008.6154
$ws.Range("M113").Value = 1314
$ws.Range("N113").Value = -6348.6154

$ws.Range("H132").Value = 2233.2273
$ws.Range("I132").Value = 1428.4117
$ws.Range("J132").Value = 4969.6
$ws.Range("K132").Value = 4285.2351
$ws.Range("L132").Value = 14908.8
$ws.Range("M132").Value = -1755.2351
$ws.Range("N132").Value = -19968.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 9984006
$ws.Range("I122").Value = 12152963
$ws.Range("K122").Value = 36458889
$ws.Range("M122").Value = -36456439

$ws.Range("H132").Value = 5566
$ws.Range("I132").Value = 5082.5
$ws.Range("J132").Value = 7500
$ws.Range("K132").Value = 15247.5
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -12717.5
$ws.Range("N132").Value = -27560

$ws.Range("H136").Value = 4832.875
$ws.Range("I136").Value = 4926.073
$ws.Range("J136").Value = 4287
$ws.Range("K136").Value = 14778.219
$ws.Range("L136").Value = 12861
$ws.Range("M136").Value = -12228.219
$ws.Range("N136").Value = -17961

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 560
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 560
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1680
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -5520

$ws.Range("H132").Value = 2273
$ws.Range("I132").Value = 1869
$ws.Range("K132").Value = 5607
$ws.Range("M132").Value = -3077
